$p = $ppt.ActivePresentation
$s = $p.Slides.Item(10)

# Remove the pie-chart picture ("Picture 4" / rId2) that was added to the
# "Böngészők" slide - it is no longer part of the slide contents.
$s.Shapes.Item("Picture 4").Delete()
